$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.424211666666666
$ws.Range("H2").Value = 4.272634999999999
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1320293333333333
$ws.Range("N2").Value = 0.396088
$ws.Range("O2").Value = 0.02700478969442551
$ws.Range("P2").Value = 0.02700478969442551
$ws.Range("Q2").Value = 0.1880377168755555
$ws.Range("R2").Value = 1.69233945188
$ws.Range("S2").Value = 0.02700478969442551
$ws.Range("T2").Value = 0.02700478969442551

# Row 3
$ws.Range("G3").Value = 1.424211666666666
$ws.Range("H3").Value = 4.272634999999999
$ws.Range("M3").Value = 0.5252536666666666
$ws.Range("N3").Value = 1.575761
$ws.Range("O3").Value = 0.1074334350287755
$ws.Range("P3").Value = 0.1074334350287755
$ws.Range("Q3").Value = 0.748072400026111
$ws.Range("R3").Value = 6.732651600234998
$ws.Range("S3").Value = 0.1074334350287755
$ws.Range("T3").Value = 0.1074334350287755

# Row 4
$ws.Range("G4").Value = 1.424211666666666
$ws.Range("H4").Value = 4.272634999999999
$ws.Range("M4").Value = 1.687203666666667
$ws.Range("N4").Value = 5.061611
$ws.Range("O4").Value = 0.3450943744066743
$ws.Range("P4").Value = 0.3450943744066743
$ws.Range("Q4").Value = 2.402935146109444
$ws.Range("R4").Value = 21.626416314985
$ws.Range("S4").Value = 0.3450943744066743
$ws.Range("T4").Value = 0.3450943744066743

# Row 5
$ws.Range("G5").Value = 1.424211666666666
$ws.Range("H5").Value = 4.272634999999999
$ws.Range("M5").Value = 0.505857
$ws.Range("N5").Value = 1.517571
$ws.Range("O5").Value = 0.1034661128369428
$ws.Range("P5").Value = 0.1034661128369428
$ws.Range("Q5").Value = 0.7204474410649999
$ws.Range("R5").Value = 6.484026969584999
$ws.Range("S5").Value = 0.1034661128369428
$ws.Range("T5").Value = 0.1034661128369428

# Row 6
$ws.Range("G6").Value = 1.424211666666666
$ws.Range("H6").Value = 4.272634999999999
$ws.Range("M6").Value = 2.038764333333333
$ws.Range("N6").Value = 6.116293000000001
$ws.Range("O6").Value = 0.4170012880331818
$ws.Range("P6").Value = 0.4170012880331818
$ws.Range("Q6").Value = 2.903631949117222
$ws.Range("R6").Value = 26.132687542055
$ws.Range("S6").Value = 0.4170012880331818
$ws.Range("T6").Value = 0.4170012880331818
